$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.525.45'
$ws.Range("E2").Value = '  +2.15%  '

$ws.Range("D3").Value = '1.788.73'
$ws.Range("E3").Value = '  +4.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5336'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +11.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3766'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.97'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07480'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.36%  '

$ws.Range("E11").Value = '  +6.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9995'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("E13").Value = '  +5.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.142'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.37%  '

$ws.Range("D15").Value = '1.784.62'
$ws.Range("E15").Value = '  +4.12%  '

$ws.Range("E16").Value = '  +3.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06442'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '

$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.929'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.74%  '

$ws.Range("D23").Value = '27.557.63'
$ws.Range("E23").Value = '  +2.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.093'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.57%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.391'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.31%  '

$ws.Range("D29").Value = '1.989.99'
$ws.Range("E29").Value = '  +4.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.73%  '

$ws.Range("E31").Value = '  +7.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1025'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.644'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.628'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.91%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02269'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.604'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +15.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06007'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.09%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.958'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.11%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2074'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6190'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.417'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9986'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.145'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.43'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5831'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.632'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.908'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.127'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06740'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.29%  '
